$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new data row (row 10) continuing the Google Code Jam statistics table:
# YEAR | ROUND | SOLVED | TRIED
$ws.Cells.Item(10, 1).Value = 2013
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = "A,B"
$ws.Cells.Item(10, 4).Value = "C"

$ws.Range("C11").Select()
